$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update B-column labels for rows 10-19 (new ordering of schemes) ---
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"

# --- Update C:P numeric data for rows 10-19 ---
$row10 = New-Object 'object[,]' 1,14
$row10[0,0] = 0.9971603071729759
$row10[0,1] = 1.012708961054975
$row10[0,2] = 0.9952648695185737
$row10[0,3] = 1.001633974110304
$row10[0,4] = 0.9971603071729759
$row10[0,5] = 1.012708961054975
$row10[0,6] = 0.9979085223676277
$row10[0,7] = 1.007545575179965
$row10[0,8] = 0.998206045989162
$row10[0,9] = 1.011503242338255
$row10[0,10] = 0.9971603071729759
$row10[0,11] = 1.003986915286774
$row10[0,12] = 1.001692027964207
$row10[0,13] = 1.00274143721648
$ws.Range("C10:P10").Value = $row10

$row11 = New-Object 'object[,]' 1,14
$row11[0,0] = 0.9896284145387224
$row11[0,1] = 1.019948022188577
$row11[0,2] = 0.9985616832543092
$row11[0,3] = 1.005350645532892
$row11[0,4] = 0.9896284145387224
$row11[0,5] = 1.019948022188577
$row11[0,6] = 0.9950939263028509
$row11[0,7] = 1.004784729535895
$row11[0,8] = 0.9963862448528026
$row11[0,9] = 1.013890872157889
$row11[0,10] = 0.9896284145387224
$row11[0,11] = 1.009254852721443
$row11[0,12] = 1.003372191378625
$row11[0,13] = 1.002955567295492
$ws.Range("C11:P11").Value = $row11

$row12 = New-Object 'object[,]' 1,14
$row12[0,0] = 0.9896111749712236
$row12[0,1] = 1.01984684633694
$row12[0,2] = 0.9986145700020101
$row12[0,3] = 1.005341454533783
$row12[0,4] = 0.9896111749712236
$row12[0,5] = 1.01984684633694
$row12[0,6] = 0.9951318497474879
$row12[0,7] = 1.004768445289851
$row12[0,8] = 0.9964015801152789
$row12[0,9] = 1.013820754057634
$row12[0,10] = 0.9896111749712236
$row12[0,11] = 1.009230708169475
$row12[0,12] = 1.00335351146099
$row12[0,13] = 1.002942084381776
$ws.Range("C12:P12").Value = $row12

$row13 = New-Object 'object[,]' 1,14
$row13[0,0] = 0.9895821848785107
$row13[0,1] = 1.01993097531062
$row13[0,2] = 0.9985749915302462
$row13[0,3] = 1.005351644924269
$row13[0,4] = 0.9895821848785107
$row13[0,5] = 1.01993097531062
$row13[0,6] = 0.9950990280189881
$row13[0,7] = 1.00477864479424
$row13[0,8] = 0.9963878145234335
$row13[0,9] = 1.013875740738752
$row13[0,10] = 0.9895821848785107
$row13[0,11] = 1.009252983420433
$row13[0,12] = 1.003359949160912
$row13[0,13] = 1.002947628089882
$ws.Range("C13:P13").Value = $row13

$row14 = New-Object 'object[,]' 1,14
$row14[0,0] = 0.9910439999999995
$row14[0,1] = 1.037095999999997
$row14[0,2] = 0.9889560000000007
$row14[0,3] = 1.006671999999999
$row14[0,4] = 0.9910439999999995
$row14[0,5] = 1.037095999999997
$row14[0,6] = 0.9891599999999997
$row14[0,7] = 1.007292
$row14[0,8] = 0.9935399999999994
$row14[0,9] = 1.025219999999998
$row14[0,10] = 0.9910439999999995
$row14[0,11] = 1.013025999999999
$row14[0,12] = 1.005941999999999
$row14[0,13] = 1.004872499999999
$ws.Range("C14:P14").Value = $row14

$row15 = New-Object 'object[,]' 1,14
$row15[0,0] = 0.99
$row15[0,1] = 1.06
$row15[0,2] = 0.98
$row15[0,3] = 1.01
$row15[0,4] = 0.99
$row15[0,5] = 1.06
$row15[0,6] = 0.98
$row15[0,7] = 1.01
$row15[0,8] = 0.99
$row15[0,9] = 1.04
$row15[0,10] = 0.99
$row15[0,11] = 1.02
$row15[0,12] = 1.01
$row15[0,13] = 1.0075
$ws.Range("C15:P15").Value = $row15

$row16 = New-Object 'object[,]' 1,14
$row16[0,0] = 0.9938177019904003
$row16[0,1] = 1.035388963635199
$row16[0,2] = 0.9894598791168016
$row16[0,3] = 1.006415851315201
$row16[0,4] = 0.9938177019904003
$row16[0,5] = 1.035388963635199
$row16[0,6] = 0.9889810655232004
$row16[0,7] = 1.006436758528
$row16[0,8] = 0.9945915418623977
$row16[0,9] = 1.023883991039998
$row16[0,10] = 0.9938177019904003
$row16[0,11] = 1.012424421376
$row16[0,12] = 1.006270599014401
$row16[0,13] = 1.0048719691264
$ws.Range("C16:P16").Value = $row16

$row17 = New-Object 'object[,]' 1,14
$row17[0,0] = 1.00232902490542
$row17[0,1] = 1.001035461848416
$row17[0,2] = 1.002308611687737
$row17[0,3] = 1.001142918343184
$row17[0,4] = 1.00232902490542
$row17[0,5] = 1.001035461848416
$row17[0,6] = 1.001434463330522
$row17[0,7] = 1.001240015309234
$row17[0,8] = 1.001124381176629
$row17[0,9] = 1.002032469518527
$row17[0,10] = 1.00232902490542
$row17[0,11] = 1.001672036768077
$row17[0,12] = 1.00170400419619
$row17[0,13] = 1.001580918264959
$ws.Range("C17:P17").Value = $row17

$row18 = New-Object 'object[,]' 1,14
$row18[0,0] = 1.006781695757998
$row18[0,1] = 0.9987070879496249
$row18[0,2] = 1.002049613690114
$row18[0,3] = 1.000803995278833
$row18[0,4] = 1.006781695757998
$row18[0,5] = 0.9987070879496249
$row18[0,6] = 1.002188618541372
$row18[0,7] = 1.000958362836331
$row18[0,8] = 1.002318483057291
$row18[0,9] = 1.001574166590589
$row18[0,10] = 1.006781695757998
$row18[0,11] = 1.00037835081987
$row18[0,12] = 1.002085598169142
$row18[0,13] = 1.001922752962769
$ws.Range("C18:P18").Value = $row18

$row19 = New-Object 'object[,]' 1,14
$row19[0,0] = 1.006786099779913
$row19[0,1] = 0.9947732152786236
$row19[0,2] = 1.003899590048737
$row19[0,3] = 0.9999226653684493
$row19[0,4] = 1.006786099779913
$row19[0,5] = 0.9947732152786236
$row19[0,6] = 1.003650996721836
$row19[0,7] = 0.9999887974060006
$row19[0,8] = 1.002610041817708
$row19[0,9] = 0.9978517162521725
$row19[0,10] = 1.006786099779913
$row19[0,11] = 0.9993364026636804
$row19[0,12] = 1.001345392618931
$row19[0,13] = 1.00118539033418
$ws.Range("C19:P19").Value = $row19

# --- A17:A19 need values + style (bold/border/center like rest of column A) ---
$ws.Range("A17").Value = 15
$ws.Range("A18").Value = 16
$ws.Range("A19").Value = 17
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null
